{"js": "// Convert the three M2Doc conditional fields (m:if / m:elseif / m:endif),\n// which are currently stored as real Word fields (fldChar begin/instrText/\n// fldChar end), into plain literal text runs wrapped in curly braces, e.g.\n// \"{m:if self.name <> 'anydsl'}\". This mirrors the TokenIteratorFieldRewriterSplit\n// rewrite: the field machinery disappears and the field code becomes\n// ordinary paragraph text, split across the same number of runs the\n// instrText was split across (plus the leading \"{\" / trailing \"}\").\n\nfunction escapeXmlText(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a minimal single-paragraph OOXML payload (wrapped in the package\n// format insertOoxml expects) containing one <w:r><w:t>.../> per entry in\n// `runs`. Each entry is [text, preserveSpace].\nfunction buildParagraphRunsOoxml(runs) {\n  let runsXml = \"\";\n  for (const [text, preserve] of runs) {\n    const esc = escapeXmlText(text);\n    const spaceAttr = preserve ? ' xml:space=\"preserve\"' : \"\";\n    runsXml += `<w:r><w:t${spaceAttr}>${esc}</w:t></w:r>`;\n  }\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body><w:p>${runsXml}</w:p></w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\n// Replace the (single) Word field living in `paragraph` with literal text\n// runs described by `runs` ([text, preserveSpace] pairs). The field\n// (fldChar begin/instrText.../fldChar end) is deleted first so no leftover\n// empty run remains, then the new runs are inserted at the (now empty)\n// paragraph start.\nasync function replaceFieldWithLiteralRuns(paragraph, runs) {\n  const fields = paragraph.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  if (fields.items.length === 0) {\n    throw new Error(\"Expected a field in this paragraph but found none.\");\n  }\n\n  fields.items[0].delete();\n  await context.sync();\n\n  paragraph.insertOoxml(buildParagraphRunsOoxml(runs), Word.InsertLocation.start);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph layout in the document:\n// 0: \"Basic if demonstration :\"\n// 1: the {m:if ...} field               <-- rewrite\n// 2: \"The THEN paragraph.\"\n// 3: the {m:elseif ...} field           <-- rewrite\n// 4: \"The ELSEIF paragraph.\"\n// 5: the {m:endif} field                <-- rewrite\n// 6: \"End of demonstration.\"\n// 7: trailing empty paragraph\nconst ifParagraph = paragraphs.items[1];\nconst elseifParagraph = paragraphs.items[3];\nconst endifParagraph = paragraphs.items[5];\n\n// {m:if self.name <> 'anydsl'}\nawait replaceFieldWithLiteralRuns(ifParagraph, [\n  [\"{m:if \", true],\n  [\"self.name \", true],\n  [\"<>\", false],\n  [\" \", true],\n  [\"'\", false],\n  [\"anydsl\", false],\n  [\"'}\", false],\n]);\n\n// {m:elseif self.name = 'anydsl'}\nawait replaceFieldWithLiteralRuns(elseifParagraph, [\n  [\"{m:\", false],\n  [\"elseif self.name = 'anydsl'\", false],\n  [\"}\", true],\n]);\n\n// {m:endif}\nawait replaceFieldWithLiteralRuns(endifParagraph, [[\"{m:endif}\", true]]);\n", "ps1": "# Convert the three M2Doc conditional fields (m:if / m:elseif / m:endif),\n# which are currently stored as real Word fields (fldChar begin/instrText/\n# fldChar end), into plain literal text runs wrapped in curly braces, e.g.\n# \"{m:if self.name <> 'anydsl'}\". This mirrors the TokenIteratorFieldRewriterSplit\n# rewrite: the field machinery disappears and the field code becomes ordinary\n# paragraph text, split across the same number of runs the instrText was\n# split across (plus the leading \"{\" / trailing \"}\").\n\nfunction Get-ParagraphAttrsAndPPr($paragraph) {\n    # Pull the paragraph's own OOXML so we can keep its existing <w:p ...>\n    # attributes (w:rsidP/w:rsidR/w:rsidRDefault) and <w:pPr> (tab stops,\n    # etc.) untouched -- only the runs inside it are being replaced.\n    $xml = $paragraph.Range.WordOpenXML\n\n    $attrs = \"\"\n    if ($xml -match '(?s)<w:p ([^>]*)>') {\n        $rawAttrs = $matches[1]\n        # WordOpenXML round-tripping stamps synthetic w14:paraId/w14:textId\n        # attributes that are not present in the source document -- keep\n        # only the real w:rsid* attributes that already existed.\n        $keep = [regex]::Matches($rawAttrs, 'w:rsid\\w*=\"[^\"]*\"') | ForEach-Object { $_.Value }\n        $attrs = [string]::Join(\" \", $keep)\n    }\n\n    $pPr = \"\"\n    if ($xml -match '(?s)<w:pPr>.*?</w:pPr>') {\n        $pPr = $matches[0]\n    }\n\n    return @{ Attrs = $attrs; PPr = $pPr }\n}\n\nfunction Escape-XmlText($s) {\n    return $s -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n}\n\n# Replace the single Word field living in $paragraph with literal text runs\n# described by $runs (an array of [text, preserveSpace] pairs). The field\n# (fldChar begin/instrText.../fldChar end) is deleted first so no leftover\n# empty run remains, then the new runs are inserted as a full paragraph\n# replacement that keeps the paragraph's own attributes/pPr.\nfunction Replace-FieldWithLiteralRuns($paragraph, $runs) {\n    $info = Get-ParagraphAttrsAndPPr $paragraph\n\n    $fields = $paragraph.Range.Fields\n    if ($fields.Count -eq 0) {\n        throw \"Expected a field in this paragraph but found none.\"\n    }\n    [void]$fields.Item(1).Delete()\n\n    $runsXml = \"\"\n    foreach ($pair in $runs) {\n        $text = $pair[0]\n        $preserve = $pair[1]\n        $esc = Escape-XmlText $text\n        if ($preserve) {\n            $runsXml += \"<w:r><w:t xml:space=`\"preserve`\">$esc</w:t></w:r>\"\n        }\n        else {\n            $runsXml += \"<w:r><w:t>$esc</w:t></w:r>\"\n        }\n    }\n\n    $pOpen = \"<w:p\"\n    if ($info.Attrs -ne \"\") { $pOpen += \" \" + $info.Attrs }\n    $pOpen += \">\"\n\n    $newXml = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n              '<w:body>' + $pOpen + $info.PPr + $runsXml + '</w:p></w:body></w:document>'\n\n    $r = $paragraph.Range\n    $r.Collapse(1) # wdCollapseStart\n    $r.InsertXML($newXml)\n}\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n\n# Paragraph layout in the document (1-based, Word COM style):\n# 1: \"Basic if demonstration :\"\n# 2: the {m:if ...} field               <-- rewrite\n# 3: \"The THEN paragraph.\"\n# 4: the {m:elseif ...} field           <-- rewrite\n# 5: \"The ELSEIF paragraph.\"\n# 6: the {m:endif} field                <-- rewrite\n# 7: \"End of demonstration.\"\n# 8: trailing empty paragraph\n$ifParagraph = $paragraphs.Item(2)\n$elseifParagraph = $paragraphs.Item(4)\n$endifParagraph = $paragraphs.Item(6)\n\n# {m:if self.name <> 'anydsl'}\nReplace-FieldWithLiteralRuns $ifParagraph @(\n    , @(\"{m:if \", $true)\n    , @(\"self.name \", $true)\n    , @(\"<>\", $false)\n    , @(\" \", $true)\n    , @(\"'\", $false)\n    , @(\"anydsl\", $false)\n    , @(\"'}\", $false)\n)\n\n# {m:elseif self.name = 'anydsl'}\nReplace-FieldWithLiteralRuns $elseifParagraph @(\n    , @(\"{m:\", $false)\n    , @(\"elseif self.name = 'anydsl'\", $false)\n    , @(\"}\", $true)\n)\n\n# {m:endif}\nReplace-FieldWithLiteralRuns $endifParagraph @(\n    , @(\"{m:endif}\", $true)\n)\n"}
